$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Build a minimal OOXML fragment (wrapped in pkg:package) containing the
# given literal run(s) XML, suitable for Range.InsertXML.
function New-RunsXml([string[]]$runsXml) {
    $body = [string]::Join("", $runsXml)
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $body + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locate the single, unique occurrence of $findText (used purely to pinpoint
# a location unambiguously, e.g. by including following-paragraph context
# via a backtick-r) and replace the first $spanLength characters of that
# match with the given run(s) XML, splitting that span into multiple <w:r>
# elements as needed. Anything after $spanLength within the match (e.g. an
# already-separate, untouched ":" run) is left completely alone.
#
# InsertXML replaces the exact range it is called on, so we first collapse
# (empty out) the target span and then insert at that now-collapsed point;
# this keeps surrounding runs (e.g. an adjacent untouched ":" run) in their
# original order instead of being reshuffled.
function Replace-SpanWithRuns([string]$findText, [int]$spanLength, [string[]]$runsXml) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Replace-SpanWithRuns: text not found: $findText"
    }
    $start = $rng.Start
    $target = $d.Range($start, $start + $spanLength)
    $target.Text = ""
    $ins = $d.Range($start, $start)
    $ins.InsertXML((New-RunsXml $runsXml))
}

# Locate the single, unique occurrence of $findText and replace the first
# $spanLength characters of that match with plain $newText, staying within
# the existing run (no run splitting). Anything after $spanLength in the
# match is left untouched.
function Replace-SpanWithText([string]$findText, [int]$spanLength, [string]$newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Replace-SpanWithText: text not found: $findText"
    }
    $start = $rng.Start
    $target = $d.Range($start, $start + $spanLength)
    $target.Text = $newText
}

# ---------------------------------------------------------------------------
# 1) "Részben teljesített" (run) immediately followed by a separate ":" run
#    -> split into two runs: "Részben " (space preserved) + "teljesítve".
#    Disambiguated from another, unrelated "Részben teljesített és ..."
#    sentence elsewhere (which spans multiple runs) by requiring the
#    immediately following ":"; only the 19 characters of "Részben
#    teljesített" itself (not the colon) are replaced, so the existing ":"
#    run is left completely untouched.
# ---------------------------------------------------------------------------
Replace-SpanWithRuns "Részben teljesített:" 19 @(
    '<w:r><w:t xml:space="preserve">Részben </w:t></w:r>',
    '<w:r><w:t>teljesítve</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# 4) "Teljesített:" (single run including the colon), immediately followed
#    by the paragraph "Amikor már ki lett egyenlítve." -> split into two
#    runs: "Teljesítve" + ":"
# ---------------------------------------------------------------------------
Replace-SpanWithRuns "Teljesített:`rAmikor már ki lett egyenlítve." 12 @(
    '<w:r><w:t>Teljesítve</w:t></w:r>',
    '<w:r><w:t>:</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# 6) "Teljesített:" (single run including the colon), immediately followed
#    by the paragraph "Amikor már teljesítve lett." -> split into two runs:
#    "Teljesítve" + ":"
# ---------------------------------------------------------------------------
Replace-SpanWithRuns "Teljesített:`rAmikor már teljesítve lett." 12 @(
    '<w:r><w:t>Teljesítve</w:t></w:r>',
    '<w:r><w:t>:</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# 5) "Sztornózott:" (single run including the colon) -> split into two
#    runs: "Sztornózva" + ":"
# ---------------------------------------------------------------------------
Replace-SpanWithRuns "Sztornózott:" 12 @(
    '<w:r><w:t>Sztornózva</w:t></w:r>',
    '<w:r><w:t>:</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# 7) "Lemondott:" (single run including the colon) -> split into two runs:
#    "Lemondva" + ":"
#    Done before hunk 3 below, because a bare "Lemondott" search (no colon
#    requirement) would otherwise also match as a prefix of this "Lemondott:"
#    occurrence.
# ---------------------------------------------------------------------------
Replace-SpanWithRuns "Lemondott:" 10 @(
    '<w:r><w:t>Lemondva</w:t></w:r>',
    '<w:r><w:t>:</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# 2) "Teljesített" (run) immediately followed by a separate ":" run, in turn
#    followed by the paragraph starting "Amikor a belőle generált számla ki
#    lett egyenlítve" -> stays a single run, just becomes "Teljesítve". The
#    ":" run is left completely untouched.
# ---------------------------------------------------------------------------
Replace-SpanWithText "Teljesített:`rAmikor a belőle generált számla ki lett egyenlítve" 11 "Teljesítve"

# ---------------------------------------------------------------------------
# 3) "Lemondott" (single run, no colon) -> "Teljesítve"-style simple rename:
#    stays a single run, becomes "Lemondva". By this point the only
#    remaining occurrence of bare "Lemondott" is this one (the "Lemondott:"
#    occurrence from hunk 7 has already been converted above).
# ---------------------------------------------------------------------------
Replace-SpanWithText "Lemondott" 9 "Lemondva"
